# Update the "想去人数" (F column) counts on the "展览" sheet and the
# "全部类型" sheet for the events whose attendance numbers changed.
#
# Each pair below is: (row on "展览" sheet, row on "全部类型" sheet, new value)
$updates = @(
    @{ Sheet1Row = 2;  Sheet4Row = 3;  NewValue = 1004 },
    @{ Sheet1Row = 6;  Sheet4Row = 10; NewValue = 7262 },
    @{ Sheet1Row = 10; Sheet4Row = 14; NewValue = 767 },
    @{ Sheet1Row = 12; Sheet4Row = 16; NewValue = 67 },
    @{ Sheet1Row = 15; Sheet4Row = 19; NewValue = 839 },
    @{ Sheet1Row = 16; Sheet4Row = 21; NewValue = 2897 },
    @{ Sheet1Row = 17; Sheet4Row = 22; NewValue = 156 },
    @{ Sheet1Row = 18; Sheet4Row = 23; NewValue = 46 },
    @{ Sheet1Row = 21; Sheet4Row = 28; NewValue = 38 },
    @{ Sheet1Row = 22; Sheet4Row = 29; NewValue = 431 },
    @{ Sheet1Row = 24; Sheet4Row = 31; NewValue = 151 },
    @{ Sheet1Row = 25; Sheet4Row = 32; NewValue = 191 },
    @{ Sheet1Row = 26; Sheet4Row = 33; NewValue = 146 },
    @{ Sheet1Row = 27; Sheet4Row = 34; NewValue = 198 },
    @{ Sheet1Row = 29; Sheet4Row = 36; NewValue = 71 },
    @{ Sheet1Row = 30; Sheet4Row = 37; NewValue = 186 },
    @{ Sheet1Row = 33; Sheet4Row = 40; NewValue = 308 },
    @{ Sheet1Row = 34; Sheet4Row = 41; NewValue = 388 },
    @{ Sheet1Row = 36; Sheet4Row = 43; NewValue = 18 },
    @{ Sheet1Row = 38; Sheet4Row = 45; NewValue = 45 }
)

$wb = $excel.ActiveWorkbook
$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

foreach ($u in $updates) {
    $wsExhibition.Cells.Item($u.Sheet1Row, 6).Value = $u.NewValue
    $wsAllTypes.Cells.Item($u.Sheet4Row, 6).Value   = $u.NewValue
}
